$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 9 (the "Lamp (k = 3, eps = -0.128)" row) into a new row 10,
# carrying over all values/number formats/styles, then tweak the case name
# and the eps (column L) value for the new scenario.
$src = $ws.Range("A9:AB9")
$dst = $ws.Range("A10:AB10")
$src.Copy($dst)

$ws.Range("B10").Value = "Lamp (k = 3, eps = -0.6)"
$ws.Range("L10").Value = -0.6

# Re-create the cell comments (Matthew Heun's engineering-unit annotations)
# on row 10 to match the ones already present on row 9.
$ws.Range("H10").AddComment("Matthew Heun:" + [char]10 + "MJ/kW-hr")
$ws.Range("I10").AddComment("Matthew Heun:" + [char]10 + "MJ/$")
$ws.Range("K10").AddComment("Matthew Heun:" + [char]10 + "$/kW-hr")
$ws.Range("O10").AddComment("Matthew Heun:" + [char]10 + "lm-hr/kW-hr")
$ws.Range("P10").AddComment("Matthew Heun:" + [char]10 + "lm-hr/kW-hr")
$ws.Range("Q10").AddComment("Matthew Heun:" + [char]10 + "Lm-hr/yr")
$ws.Range("R10").AddComment("Matthew Heun:" + [char]10 + "$/year")
$ws.Range("S10").AddComment("Matthew Heun:" + [char]10 + "$")
$ws.Range("T10").AddComment("Matthew Heun:" + [char]10 + "$")
$ws.Range("U10").AddComment("Matthew Heun:" + [char]10 + "$/year")
$ws.Range("V10").AddComment("Matthew Heun:" + [char]10 + "$/year")
$ws.Range("Y10").AddComment("Matthew Heun:" + [char]10 + "MJ")
$ws.Range("Z10").AddComment("Matthew Heun:" + [char]10 + "years")
$ws.Range("AA10").AddComment("Matthew Heun:" + [char]10 + "MJ")
$ws.Range("AB10").AddComment("Matthew Heun:" + [char]10 + "years")

# Update the view's active selection, matching where the author left off.
$ws.Range("B7").Select() | Out-Null
